$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: new week of data (date shifts from 23-Dec-2024 to 13-Jan-2025, and several
# task columns get a score of 5)
$ws.Range("A8").Value = 45670
$ws.Range("B8").Value = 5
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 5
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 5
$ws.Range("N8").Value = 5
$ws.Range("R8").Value = 5
$ws.Range("T8").Value = 5
$ws.Range("X8").Value = 5
$ws.Range("Y8").Value = 5
$ws.Range("Z8").Value = 5

# Rows 9-15: shift the week date forward by 21 days (3 weeks) for each remaining row
$ws.Range("A9").Value = 45677
$ws.Range("A10").Value = 45684
$ws.Range("A11").Value = 45691
$ws.Range("A12").Value = 45698
$ws.Range("A13").Value = 45705
$ws.Range("A14").Value = 45712
$ws.Range("A15").Value = 45719

# Update the active selection to reflect where the user left off editing
$ws.Range("R17").Select() | Out-Null
